$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.959.97"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "3.567.26"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.09"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.50"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").Value = "3.565.59"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "4.175.18"
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000182"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "3.569.81"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.05"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "64.579.07"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.38"
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.83"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.98"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  +4.39%  "
$ws.Range("D24").Value = "3.713.31"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.12"
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +5.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.73"
$ws.Range("E28").Value = "  +5.64%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +3.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.43"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("E32").Value = "  +23.81%  "
$ws.Range("D33").Value = "3.569.30"
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.01"
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.94"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "169.33"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("E40").Value = "  +5.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0809"
$ws.Range("E41").Value = "  +3.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.20"
$ws.Range("E42").Value = "  +8.81%  "
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.73"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("E47").Value = "  +4.65%  "
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").Value = "2.481.23"
$ws.Range("E49").Value = "  +11.65%  "
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.37"
$ws.Range("E51").Value = "  +8.79%  "
